# Generate Report for Archive
# - Status text "Ready for handoff" -> "In Translation" on the Overview sheet
#   (zh-cn/de-de status columns) and on each per-locale sheet's Status column.
# - The Status columns are narrower now that "In Translation" is shorter than
#   "Ready for handoff", so their column widths shrink to match.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: E2 (zh-cn status) and F2 (de-de status) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# Narrow columns E and F to the new best-fit width.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet: C2 (Status) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet: C2 (Status) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
